$wb = $excel.ActiveWorkbook

# "bookings.views.py" is the 3rd sheet in the workbook.
$ws = $wb.Worksheets.Item(3)

# The old B18 ("...prevents duplicate bookings") is replaced with a new
# test description about preventing booking onto an already-booked slot
# that isn't the one being edited.
$ws.Range("B18").Value2 = "Test that the edit booking POST method prevents booking onto a booked slot that is not the slot being edited"

# Insert a new row right after it for another new test case, pushing the
# remaining rows (old B19 onward) down by one.
[void]$ws.Rows.Item(19).Insert()
$ws.Range("B19").Value2 = "Tests that the edit booking POST method allows a user to edit just the lesson type on an existing booking"

# Update the sheet's selection to match the saved view, then make this
# sheet the active tab of the workbook.
[void]$ws.Range("B21").Select()
$ws.Activate()
